# Update bus voltage magnitude results (vm_pu) for the 380 kV case:
# bus-0 (slack) setpoint drops from 1.05 pu to 1.02 pu, and the
# downstream bus voltages (columns B-E, I-L, N) are recomputed accordingly
# for every result row (r2:r25) while columns A and G are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    "B2"=1.02; "C2"=1.017623197690867; "D2"=1.02279594305572; "E2"=1.018992132068612; "I2"=1.026719489015036; "J2"=1.022836260980906; "K2"=1.025629270358936; "L2"=1.021836702599545; "N2"=1.011683407944881
    "B3"=1.02; "C3"=1.018647958511197; "D3"=1.023509632671733; "E3"=1.019863083122803; "I3"=1.026847841497977; "J3"=1.023496594028466; "K3"=1.026149956389011; "L3"=1.02251340393172; "N3"=1.011904588546823
    "B4"=1.02; "C4"=1.019311162922379; "D4"=1.023971288540365; "E4"=1.020427135149274; "I4"=1.026929445470069; "J4"=1.023923447514811; "K4"=1.026486058119088; "L4"=1.02295113062842; "N4"=1.012047480992211
    "B5"=1.02; "C5"=1.019590001521214; "D5"=1.024165331929493; "E5"=1.020664378719889; "I5"=1.026963404544207; "J5"=1.024102794105646; "K5"=1.026627158974945; "L5"=1.023135115737427; "N5"=1.012107498565619
    "B6"=1.02; "C6"=1.019636821353119; "D6"=1.024197910483183; "E6"=1.020704219743003; "I6"=1.026969086047836; "J6"=1.02413290115276; "K6"=1.026650838896389; "L6"=1.023166005564224; "N6"=1.012117572582065
    "B7"=1.02; "C7"=1.019314888669737; "D7"=1.023973881503321; "E7"=1.020430304754342; "I7"=1.026929900597808; "J7"=1.02392584435669; "K7"=1.026487944287926; "L7"=1.022953589185553; "N7"=1.012048283164087
    "B8"=1.02; "C8"=1.017969496332596; "D8"=1.023037168026617; "E8"=1.019286372413637; "I8"=1.026763166040277; "J8"=1.023059511506999; "K8"=1.02580540735578; "L8"=1.022065426477186; "N8"=1.011758203644242
    "B9"=1.02; "C9"=1.015599638031913; "D9"=1.021385464178808; "E9"=1.017274395178065; "I9"=1.026458281903396; "J9"=1.021529684422721; "K9"=1.024596465275399; "L9"=1.020499293212192; "N9"=1.011245325419931
    "B10"=1.02; "C10"=1.014020343832037; "D10"=1.020283653267004; "E10"=1.01593566474768; "I10"=1.026247602974655; "J10"=1.020507655945194; "K10"=1.023786360592629; "L10"=1.019454516635169; "N10"=1.010902264179282
    "B11"=1.02; "C11"=1.01333663812325; "D11"=1.019806409347626; "E11"=1.015356603006007; "I11"=1.026154621024324; "J11"=1.020064604325109; "K11"=1.023434601997258; "L11"=1.019001960902855; "N11"=1.010753446826148
    "B12"=1.02; "C12"=1.013082700256922; "D12"=1.019629117744589; "E12"=1.015141607153567; "I12"=1.026119819840696; "J12"=1.019899959574638; "K12"=1.023303796889709; "L12"=1.018833837944998; "N12"=1.01069812914622
    "B13"=1.02; "C13"=1.013137169846926; "D13"=1.019667148373844; "E13"=1.015187720253068; "I13"=1.026127296733019; "J13"=1.019935279848697; "K13"=1.023331861639927; "L13"=1.018869901950237; "N13"=1.01070999679612
    "B14"=1.02; "C14"=1.013315647114262; "D14"=1.019791754810936; "E14"=1.015338829471106; "I14"=1.026151749721834; "J14"=1.020050996279346; "K14"=1.023423792587234; "L14"=1.018988064285342; "N14"=1.010748875070131
    "B15"=1.02; "C15"=1.013425615614791; "D15"=1.019868526067606; "E15"=1.015431945253245; "I15"=1.026166781103477; "J15"=1.020122282980251; "K15"=1.02348041490985; "L15"=1.01906086488688; "N15"=1.010772823925532
    "B16"=1.02; "C16"=1.014065722056937; "D16"=1.020315323242093; "E16"=1.015974108236145; "I16"=1.02625373690197; "J16"=1.020537049188243; "K16"=1.02380968510127; "L16"=1.019484547929548; "N16"=1.010912135035372
    "B17"=1.02; "C17"=1.014467281141817; "D17"=1.020595547144491; "E17"=1.016314358539546; "I17"=1.026307811943675; "J17"=1.020797085774751; "K17"=1.024015966213127; "L17"=1.019750270481991; "N17"=1.010999449191395
    "B18"=1.02; "C18"=1.014701517184884; "D18"=1.020758982225734; "E18"=1.016512880304692; "I18"=1.026339183491654; "J18"=1.020948711811783; "K18"=1.024136192099032; "L18"=1.019905246422964; "N18"=1.011050352037453
    "B19"=1.02; "C19"=1.014781387899801; "D19"=1.020814706828903; "E19"=1.016580581170345; "I19"=1.026349851616192; "J19"=1.021000404063277; "K19"=1.024177169989054; "L19"=1.019958086564313; "N19"=1.011067704171964
    "B20"=1.02; "C20"=1.014424196263998; "D20"=1.020565483293599; "E20"=1.016277846740492; "I20"=1.026302027733948; "J20"=1.020769191365063; "K20"=1.02399384396141; "L20"=1.019721762569243; "N20"=1.010990083899011
    "B21"=1.02; "C21"=1.013263089431352; "D21"=1.019755061919573; "E21"=1.015294328987025; "I21"=1.026144556198411; "J21"=1.020016922759423; "K21"=1.023396725248296; "L21"=1.018953269070349; "N21"=1.010737427497928
    "B22"=1.02; "C22"=1.012533176253563; "D22"=1.019245391239228; "E22"=1.014676493426588; "I22"=1.02604402296861; "J22"=1.019543504153504; "K22"=1.023020446853417; "L22"=1.018469950367561; "N22"=1.010578339395592
    "B23"=1.02; "C23"=1.012920105598738; "D23"=1.019515588992959; "E23"=1.015003968222141; "I23"=1.026097461932783; "J23"=1.019794513631347; "K23"=1.023219999150245; "L23"=1.018726179507829; "N23"=1.010662697023618
    "B24"=1.02; "C24"=1.014443664424098; "D24"=1.020579067899451; "E24"=1.016294344666953; "I24"=1.026304641893112; "J24"=1.020781795799176; "K24"=1.024003840346395; "L24"=1.01973464411587; "N24"=1.010994315751924
    "B25"=1.02; "C25"=1.016212196167205; "D25"=1.021812592656706; "E25"=1.017794086908304; "I25"=1.026719489015036; "J25"=1.021925561672468; "K25"=1.024909739124282; "L25"=1.020904300080484; "N25"=1.011378119167733
}

foreach ($addr in $newValues.Keys) {
    $ws.Range($addr).Value = $newValues[$addr]
}
